# Generate Report for Handback
# Adds a new handback row (bb60392b-d523-4988-a65f-e5e19e030062) to the
# Overview / zh-cn / de-de sheets, mirroring the existing
# 9c58bbe0-2c9c-499e-af8b-e0a1eefd73ba row.

$wb = $excel.ActiveWorkbook

$newFile   = "bb60392b-d523-4988-a65f-e5e19e030062.md"
$newPath   = "e2e\bb60392b-d523-4988-a65f-e5e19e030062.md"
$statusInSync = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Sheet "Overview" -> new row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add()

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("B3").Value = $newPath
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = $statusInSync
$wsOverview.Range("F3").Value = $statusInSync
$wsOverview.Range("G3").Value = "2016-10-13 13:03:33"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c2c1ce77e5151593b7a393f9291131b8ac78507/e2e/bb60392b-d523-4988-a65f-e5e19e030062.md", "", "", $newPath)

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> new row 3
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add()

$wsZhCn.Range("A3").Value = $newFile
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $statusInSync
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'True"
$wsZhCn.Range("G3").Value = "bb60392b-d523-4988-a65f-e5e19e030062.d82d13e929653a7ea77376a2d565c216de2c05a0.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-10-13 13:03:22"
$wsZhCn.Range("I3").Value = $newFile
$wsZhCn.Range("J3").Value = "bb60392b-d523-4988-a65f-e5e19e030062.d82d13e929653a7ea77376a2d565c216de2c05a0.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-10-13 13:04:08"
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("O3").Value = "'False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/09a163155d96ae2edd3d5ec61c4c23ffcf73a238/e2e/bb60392b-d523-4988-a65f-e5e19e030062.md", "", "", $newFile)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/09a163155d96ae2edd3d5ec61c4c23ffcf73a238/e2e/bb60392b-d523-4988-a65f-e5e19e030062.md", "", "", $newFile)

# ---------------------------------------------------------------------------
# Sheet "de-de" -> new row 3
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add()

$wsDeDe.Range("A3").Value = $newFile
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $statusInSync
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'True"
$wsDeDe.Range("G3").Value = "bb60392b-d523-4988-a65f-e5e19e030062.d82d13e929653a7ea77376a2d565c216de2c05a0.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-10-13 13:03:33"
$wsDeDe.Range("I3").Value = $newFile
$wsDeDe.Range("J3").Value = "bb60392b-d523-4988-a65f-e5e19e030062.d82d13e929653a7ea77376a2d565c216de2c05a0.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-10-13 13:04:25"
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("O3").Value = "'False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/447f0f04fa48cc90433dc9dbb05cdba049de5726/e2e/bb60392b-d523-4988-a65f-e5e19e030062.md", "", "", $newFile)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/447f0f04fa48cc90433dc9dbb05cdba049de5726/e2e/bb60392b-d523-4988-a65f-e5e19e030062.md", "", "", $newFile)

# ---------------------------------------------------------------------------
# Styling - mirror the existing row 2 look onto row 3
# (datetime columns keep the "yyyy-mm-dd HH:mm:ss" numeric format; the
# hyperlink columns already pick up a hyperlink style from Hyperlinks.Add)
# ---------------------------------------------------------------------------
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
